$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Execute already Y; update ActualResult/Result to reflect a Pass
$ws.Range("F2").Value = "เข้าสู่ระบบสำเร็จ"
$ws.Range("G2").Value = "Pass"

# Row 3: mark Execute = Y
$ws.Range("A3").Value = "Y"

# Row 4: mark Execute = Y
$ws.Range("A4").Value = "Y"

# Row 5: mark Execute = Y
$ws.Range("A5").Value = "Y"

# Row 6: mark Execute = Y
$ws.Range("A6").Value = "Y"

# Row 7: mark Execute = Y
$ws.Range("A7").Value = "Y"

# Row 8: mark Execute = Y, fix email test case, password, actual result, result
$ws.Range("A8").Value = "Y"
$ws.Range("C8").Value = "mju6204106340@mju.ac.th"
$ws.Range("D8").Value = 111111
$ws.Range("F8").Value = "ท่านยังไม่ได้รับการอนุมัติ"
$ws.Range("G8").Value = "Pass"

# Row 9: mark Execute = Y, update actual result message
$ws.Range("A9").Value = "Y"
$ws.Range("F9").Value = "Please fill out this field."

# Row 10: mark Execute = Y, update actual result/result
$ws.Range("A10").Value = "Y"
$ws.Range("F10").Value = "เข้าสู่ระบบสำเร็จ"
$ws.Range("G10").Value = "Pass"

# Row 11: mark Execute = Y
$ws.Range("A11").Value = "Y"

# Row 12: mark Execute = Y, update actual result message
$ws.Range("A12").Value = "Y"
$ws.Range("F12").Value = "Please fill out this field."

# Update the active selection to H5 as recorded in the latest test run
[void]$ws.Range("H5").Select()
